$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 124 <= data from row 126 (pre-edit)
$ws.Range("B124").Value = 6703695
$ws.Range("F124").Value = "Spartak Varna"
$ws.Range("G124").Value = "Septemvri Sofia"
$ws.Range("H124").Value = 1
$ws.Range("I124").Value = 0
$ws.Range("J124").Value = "H"
$ws.Range("K124").Value = 2
$ws.Range("L124").Value = 3.3
$ws.Range("M124").Value = 3.3
$ws.Range("N124").Value = 1.85
$ws.Range("O124").Value = 3.6
$ws.Range("P124").Value = 4
$ws.Range("Q124").Value = -0.5
$ws.Range("R124").Value = 1.925
$ws.Range("S124").Value = 1.925
$ws.Range("T124").Value = 2.75
$ws.Range("U124").Value = 2
$ws.Range("V124").Value = 1.85
$ws.Range("W124").Value = 0.8500000000000001
$ws.Range("X124").Value = -1
$ws.Range("Y124").Value = -1
$ws.Range("Z124").Value = 0.925
$ws.Range("AA124").Value = -1
$ws.Range("AB124").Value = -1
$ws.Range("AC124").Value = 0.8500000000000001

# Row 125 <= data from row 124 (pre-edit)
$ws.Range("B125").Value = 6703696
$ws.Range("F125").Value = "Botev Vratsa"
$ws.Range("G125").Value = "FC Hebar Pazardzhik"
$ws.Range("H125").Value = 2
$ws.Range("I125").Value = 3
$ws.Range("J125").Value = "A"
$ws.Range("K125").Value = 3.8
$ws.Range("L125").Value = 3.4
$ws.Range("M125").Value = 1.833
$ws.Range("N125").Value = 5
$ws.Range("O125").Value = 3.6
$ws.Range("P125").Value = 1.727
$ws.Range("Q125").Value = 0.75
$ws.Range("R125").Value = 1.875
$ws.Range("S125").Value = 1.975
$ws.Range("T125").Value = 2.25
$ws.Range("U125").Value = 1.975
$ws.Range("V125").Value = 1.875
$ws.Range("W125").Value = -1
$ws.Range("X125").Value = -1
$ws.Range("Y125").Value = 0.7270000000000001
$ws.Range("Z125").Value = -0.5
$ws.Range("AA125").Value = 0.4875
$ws.Range("AB125").Value = 0.9750000000000001
$ws.Range("AC125").Value = -1

# Row 126 <= data from row 125 (pre-edit)
$ws.Range("B126").Value = 6627749
$ws.Range("F126").Value = "Beroe"
$ws.Range("G126").Value = "Pirin Blagoevgrad"
$ws.Range("H126").Value = 1
$ws.Range("I126").Value = 1
$ws.Range("J126").Value = "D"
$ws.Range("K126").Value = 1.75
$ws.Range("L126").Value = 3.3
$ws.Range("M126").Value = 4.333
$ws.Range("N126").Value = 2.45
$ws.Range("O126").Value = 2.9
$ws.Range("P126").Value = 3.3
$ws.Range("Q126").Value = -0.25
$ws.Range("R126").Value = 2.025
$ws.Range("S126").Value = 1.825
$ws.Range("T126").Value = 2
$ws.Range("U126").Value = 2.1
$ws.Range("V126").Value = 1.775
$ws.Range("W126").Value = -1
$ws.Range("X126").Value = 1.9
$ws.Range("Y126").Value = -1
$ws.Range("Z126").Value = -0.5
$ws.Range("AA126").Value = 0.4125
$ws.Range("AB126").Value = 0
$ws.Range("AC126").Value = -0

# Row 133 <= data from row 134 (pre-edit)
$ws.Range("B133").Value = 6627737
$ws.Range("F133").Value = "Slavia Sofia"
$ws.Range("G133").Value = "Lokomotiv 1929 Sofia"
$ws.Range("H133").Value = 2
$ws.Range("I133").Value = 0
$ws.Range("J133").Value = "H"
$ws.Range("K133").Value = 1.5
$ws.Range("L133").Value = 3.75
$ws.Range("M133").Value = 6.5
$ws.Range("N133").Value = 1.444
$ws.Range("O133").Value = 4.333
$ws.Range("P133").Value = 8
$ws.Range("Q133").Value = -1.25
$ws.Range("R133").Value = 2
$ws.Range("S133").Value = 1.85
$ws.Range("T133").Value = 2.25
$ws.Range("U133").Value = 1.875
$ws.Range("V133").Value = 1.975
$ws.Range("W133").Value = 0.444
$ws.Range("X133").Value = -1
$ws.Range("Y133").Value = -1
$ws.Range("Z133").Value = 1
$ws.Range("AA133").Value = -1
$ws.Range("AB133").Value = -0.5
$ws.Range("AC133").Value = 0.4875

# Row 134 <= data from row 133 (pre-edit)
$ws.Range("B134").Value = 6627736
$ws.Range("F134").Value = "Botev Plovdiv"
$ws.Range("G134").Value = "Arda Kardzhali"
$ws.Range("H134").Value = 0
$ws.Range("I134").Value = 3
$ws.Range("J134").Value = "A"
$ws.Range("K134").Value = 5.25
$ws.Range("L134").Value = 3.6
$ws.Range("M134").Value = 1.571
$ws.Range("N134").Value = 26
$ws.Range("O134").Value = 11
$ws.Range("P134").Value = 1.083
$ws.Range("Q134").Value = 2.5
$ws.Range("R134").Value = 1.825
$ws.Range("S134").Value = 2.025
$ws.Range("T134").Value = 3.25
$ws.Range("U134").Value = 2
$ws.Range("V134").Value = 1.85
$ws.Range("W134").Value = -1
$ws.Range("X134").Value = -1
$ws.Range("Y134").Value = 0.08299999999999996
$ws.Range("Z134").Value = -1
$ws.Range("AA134").Value = 1.025
$ws.Range("AB134").Value = -0.5
$ws.Range("AC134").Value = 0.425

# Row 135 <= data from row 136 (pre-edit)
$ws.Range("B135").Value = 6627725
$ws.Range("F135").Value = "Levski Sofia"
$ws.Range("G135").Value = "CSKA Sofia"
$ws.Range("H135").Value = 0
$ws.Range("I135").Value = 2
$ws.Range("J135").Value = "A"
$ws.Range("K135").Value = 2.625
$ws.Range("L135").Value = 3
$ws.Range("M135").Value = 2.6
$ws.Range("N135").Value = 2.55
$ws.Range("O135").Value = 3.3
$ws.Range("P135").Value = 2.8
$ws.Range("Q135").Value = 0
$ws.Range("R135").Value = 1.825
$ws.Range("S135").Value = 2.025
$ws.Range("T135").Value = 2.25
$ws.Range("U135").Value = 1.85
$ws.Range("V135").Value = 2
$ws.Range("W135").Value = -1
$ws.Range("X135").Value = -1
$ws.Range("Y135").Value = 1.8
$ws.Range("Z135").Value = -1
$ws.Range("AA135").Value = 1.025
$ws.Range("AB135").Value = -0.5
$ws.Range("AC135").Value = 0.5

# Row 136 <= data from row 135 (pre-edit)
$ws.Range("B136").Value = 6627724
$ws.Range("F136").Value = "CSKA 1948 Sofia"
$ws.Range("G136").Value = "Lokomotiv Plovdiv"
$ws.Range("H136").Value = 0
$ws.Range("I136").Value = 0
$ws.Range("J136").Value = "D"
$ws.Range("K136").Value = 1.5
$ws.Range("L136").Value = 3.8
$ws.Range("M136").Value = 6
$ws.Range("N136").Value = 1.45
$ws.Range("O136").Value = 4.2
$ws.Range("P136").Value = 8
$ws.Range("Q136").Value = -1.25
$ws.Range("R136").Value = 2.025
$ws.Range("S136").Value = 1.825
$ws.Range("T136").Value = 2.5
$ws.Range("U136").Value = 1.85
$ws.Range("V136").Value = 2
$ws.Range("W136").Value = -1
$ws.Range("X136").Value = 3.2
$ws.Range("Y136").Value = -1
$ws.Range("Z136").Value = -1
$ws.Range("AA136").Value = 0.825
$ws.Range("AB136").Value = -1
$ws.Range("AC136").Value = 1

# Recalculated odds for rows 306-313 (independent cell updates)
$ws.Range("N306").Value = 1.85
$ws.Range("P306").Value = 4.2
$ws.Range("U306").Value = 1.975
$ws.Range("V306").Value = 1.875
$ws.Range("N307").Value = 1.285
$ws.Range("O307").Value = 5.75
$ws.Range("P307").Value = 9
$ws.Range("R307").Value = 1.85
$ws.Range("S307").Value = 2
$ws.Range("N308").Value = 1.571
$ws.Range("O308").Value = 4
$ws.Range("P308").Value = 6
$ws.Range("R308").Value = 2
$ws.Range("S308").Value = 1.85
$ws.Range("U308").Value = 2.05
$ws.Range("V308").Value = 1.8
$ws.Range("N310").Value = 1.615
$ws.Range("P310").Value = 6
$ws.Range("N311").Value = 2.15
$ws.Range("R311").Value = 1.85
$ws.Range("S311").Value = 2
$ws.Range("O312").Value = 3.5
$ws.Range("P313").Value = 1.285
$ws.Range("Q313").Value = 1.5
$ws.Range("R313").Value = 2
$ws.Range("S313").Value = 1.85
$ws.Range("U313").Value = 1.95
$ws.Range("V313").Value = 1.9
